# Trade #103 closed at 2026-02-16 21:40:28 - leadlag UP +0.000%
#
# This script reflects two related changes that happened in the live
# trading log:
#   1. Trade #69 on the "momentum" sheet (Trade # 69 / row 17), which was
#      previously OPEN, closed as a winning trade (time_exit_5min, +5 min).
#      That trade is mirrored as a new row on the "All Trades" sheet.
#   2. A brand new trade (#103) was opened on the "leadlag" sheet.
# Both ripple into the aggregate stats on "Summary" and "Comparison".

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($cell, [string]$text)
    # Force literal text storage even when the string looks like a number,
    # a percentage, a date or a time (Excel would otherwise "smart type"
    # these into numeric/date cells). Clearing formats afterwards drops the
    # temporary text number-format so no stray style survives on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# Summary sheet: OVERALL row (2) and momentum row (4) stats refresh
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Cells.Item(2, 3).Value = 69
Set-TextValue $summary.Cells.Item(2, 4) "68.1%"
Set-TextValue $summary.Cells.Item(2, 5) "+19.2527%"
Set-TextValue $summary.Cells.Item(2, 6) "+0.2790%"

Set-TextValue $summary.Cells.Item(4, 4) "56.0%"
Set-TextValue $summary.Cells.Item(4, 5) "+8.1614%"
Set-TextValue $summary.Cells.Item(4, 6) "+0.3265%"

# ---------------------------------------------------------------------
# momentum sheet: trade #69 (row 17) goes from OPEN to CLOSED
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

$momentum.Cells.Item(17, 7).Value = 67879.417698
Set-TextValue $momentum.Cells.Item(17, 8) "CLOSED"
$momentum.Cells.Item(17, 9).Value = 0.9307
$momentum.Cells.Item(17, 10).Value = 9.31
Set-TextValue $momentum.Cells.Item(17, 13) "time_exit_5min"
$momentum.Cells.Item(17, 14).Value = 5

# ---------------------------------------------------------------------
# All Trades sheet: append the now-closed trade as new row 70
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(70, 1).Value = 69
Set-TextValue $allTrades.Cells.Item(70, 2) "2026-02-16"
Set-TextValue $allTrades.Cells.Item(70, 3) "21:35:25"
Set-TextValue $allTrades.Cells.Item(70, 4) "momentum"
Set-TextValue $allTrades.Cells.Item(70, 5) "DOWN"
$allTrades.Cells.Item(70, 6).Value = 68517.13
$allTrades.Cells.Item(70, 7).Value = 67879.417698
Set-TextValue $allTrades.Cells.Item(70, 8) "CLOSED"
$allTrades.Cells.Item(70, 9).Value = 0.9307
$allTrades.Cells.Item(70, 10).Value = 9.31
$allTrades.Cells.Item(70, 11).Value = 0.9
Set-TextValue $allTrades.Cells.Item(70, 12) "Downward momentum: -0.298% over 10 samples"
Set-TextValue $allTrades.Cells.Item(70, 13) "time_exit_5min"
$allTrades.Cells.Item(70, 14).Value = 5

# ---------------------------------------------------------------------
# leadlag sheet: append newly opened trade #103 as new row 79
# ---------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

$leadlag.Cells.Item(79, 1).Value = 103
Set-TextValue $leadlag.Cells.Item(79, 2) "2026-02-16"
Set-TextValue $leadlag.Cells.Item(79, 3) "21:40:28"
Set-TextValue $leadlag.Cells.Item(79, 4) "leadlag"
Set-TextValue $leadlag.Cells.Item(79, 5) "UP"
$leadlag.Cells.Item(79, 6).Value = 68430.595
Set-TextValue $leadlag.Cells.Item(79, 8) "OPEN"
$leadlag.Cells.Item(79, 9).Value = 0
$leadlag.Cells.Item(79, 10).Value = 0
$leadlag.Cells.Item(79, 11).Value = 0.75
Set-TextValue $leadlag.Cells.Item(79, 12) "Coinbase leading with 0.101% move"
$leadlag.Cells.Item(79, 14).Value = 0

# ---------------------------------------------------------------------
# Comparison sheet: momentum row (3) stats refresh
# ---------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

Set-TextValue $comparison.Cells.Item(3, 3) "56.0%"
Set-TextValue $comparison.Cells.Item(3, 4) "8.26"
Set-TextValue $comparison.Cells.Item(3, 5) "+0.6633%"
Set-TextValue $comparison.Cells.Item(3, 7) "1.18"
